# Common: Finished mix detail
# Adds new "lab.mixture.preview*" translation rows to the Import sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy the formatting of the last existing data row (567) down onto the
# new rows (568-577) so the new cells pick up the same "import" cell
# style (wrapText) used by every other data row.
$ws.Range("A567:C567").Copy()
$ws.Range("A568:C577").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$rows = @(
    @{ Row = 568; Key = "lab.mixture.preview.pgvg";              Translation = "PG/VG" },
    @{ Row = 569; Key = "lab.mixture.preview.age";                Translation = "Stáří mixu" },
    @{ Row = 570; Key = "lab.mixture.preview.steep";              Translation = "Doba zrání" },
    @{ Row = 571; Key = "lab.mixture.preview.mixed";              Translation = "Datum mixu" },
    @{ Row = 572; Key = "lab.mixture.preview.expires";            Translation = "Expirace" },
    @{ Row = 573; Key = "lab.mixture.preview.volume";             Translation = "Objem" },
    @{ Row = 574; Key = "lab.mixture.preview";                    Translation = "Náhled mixu" },
    @{ Row = 575; Key = "lab.mixture.preview.preview.title";      Translation = "Náhled mixu" },
    @{ Row = 576; Key = "lab.mixture.preview.preview.subtitle";   Translation = "Přehled všech dostupných dat o vybraném mixu." },
    @{ Row = 577; Key = "lab.mixture.button.index";               Translation = "Detail mixu" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value2 = "cs"
    $ws.Cells.Item($r.Row, 2).Value2 = $r.Key
    $ws.Cells.Item($r.Row, 3).Value2 = $r.Translation
}

# Keep the sheet view in sync with where the new rows were added.
$ws.Application.ActiveWindow.ScrollRow = 558
$ws.Range("B570").Select()
